$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Control Lines"
$ws3 = $wb.Worksheets.Item(3)   # "Test Programs"

# ---------------------------------------------------------------------------
# New block of rows 31-37 : "// Increments A, stops at 0" test program
# ---------------------------------------------------------------------------

# Row 31: section title, styled like the other section headers (copy format from C16 -> C31)
$ws3.Range("A2:B2").Copy()
$ws3.Range("A31:B31").PasteSpecial(-4122)   # xlPasteFormats
$ws3.Range("C16").Copy()
$ws3.Range("C31").PasteSpecial(-4122)       # xlPasteFormats
$ws3.Range("C31").Value = "// Increments A, stops at 0"

# Prepare the Dec/Hex columns + D-column text format (copy format from D3)
$ws3.Range("D3").Copy()
$ws3.Range("D32:D37").PasteSpecial(-4122)   # xlPasteFormats

$ws3.Range("A32").Value = 0
$ws3.Range("A33").Value = 3
$ws3.Range("A34").Value = 6
$ws3.Range("A35").Value = 9
$ws3.Range("B32:B35").Formula = "=DEC2HEX(A32,3)"

$ws3.Range("A36").Value = 12
$ws3.Range("A37").Value = 15
$ws3.Range("B36:B37").Formula = "=DEC2HEX(A36,3)"

# Column C (instructions) - reuse existing instructions first, then new ones
# in the exact order they were authored so shared-string indices line up.
$ws3.Range("C32").Value = "LD B, 0x1"
$ws3.Range("C33").Value = "LD A, 0x20"
$ws3.Range("C34").Value = "ADD A, B"
$ws3.Range("C36").Value = "JP [0x006]"
$ws3.Range("C37").Value = "JP [0x00f]"
$ws3.Range("C35").Value = "JP Z, [0x00f]"

# Column D (machine code)
$ws3.Range("D32").Value = "04 80 01"
$ws3.Range("D33").Value = "04 00 20"
$ws3.Range("D34").Value = "14 10 00"
$ws3.Range("D36").Value = "2c 00 06"
$ws3.Range("D37").Value = "2c 00 0f"
$ws3.Range("D35").Value = "30 00 0f"

# ---------------------------------------------------------------------------
# New block of rows 41-43 : "// Load value from memory, puts it in A" program
# ---------------------------------------------------------------------------

$ws3.Range("C16").Copy()
$ws3.Range("C41").PasteSpecial(-4122)       # xlPasteFormats
$ws3.Range("C41").Value = "// Load value from memory, puts it in A"

$ws3.Range("D3").Copy()
$ws3.Range("D42:D43").PasteSpecial(-4122)   # xlPasteFormats

$ws3.Range("A42").Value = 0
$ws3.Range("B42:B47").Formula = "=DEC2HEX(A42,3)"
$ws3.Range("A43").Value = 3

$ws3.Range("C42").Value = "LD A, [0x003]"
$ws3.Range("D43").Value = "a0"
$ws3.Range("D42").Value = "0c 00 03"

# Only rows 42 and 43 are populated so far; remove the extra shared-formula
# rows 44-47 that Excel materialized while filling the B42:B47 formula range.
$ws3.Range("B44:B47").ClearContents()

# ---------------------------------------------------------------------------
# View / selection state
# ---------------------------------------------------------------------------

# "Test Programs" scrolls back to the top and the selection moves to D42
$ws3.Range("A1").Select()
$ws3.Range("D42").Select()

# "Control Lines" becomes the active/selected tab with B18 selected
$ws1.Range("B18").Select()
